$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Enterprises density (per 1000 people)" row (row 13) and the
# "Enterprises (% of total)" row (row 15) with more precise figures.
# These cells store plain text (number-like strings), so we briefly force
# a text number format while writing the value, then restore the cell's
# default ("Normal") style so formatting stays unchanged.
$cells = @(
    @{ Addr = "B13"; Val = "9.37" },
    @{ Addr = "C13"; Val = "1.74" },
    @{ Addr = "D13"; Val = "11.11" },
    @{ Addr = "B15"; Val = "81.72" },
    @{ Addr = "C15"; Val = "15.17" },
    @{ Addr = "D15"; Val = "96.88" }
)

foreach ($cell in $cells) {
    $rng = $ws.Range($cell.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.Val
    $rng.Style = "Normal"
}
